$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.287.23"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "'2.422.90"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D5").Value = "'319.12"
$ws.Range("E5").Value = "  +3.63%  "
$ws.Range("D6").Value = "'102.94"
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  +4.92%  "
$ws.Range("D10").Value = "'35.40"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "'0.0798"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "'18.23"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("D14").Value = "'7.03"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "'2.803.76"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "'2.429.46"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "'0.844"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "'45.231.70"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'6.33"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "'244.59"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").Value = "'2.48"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'25.74"
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("D29").Value = "'9.60"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "'49.55"
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("D31").Value = "'32.87"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'20.20"
$ws.Range("E32").Value = "  +8.24%  "
$ws.Range("E33").Value = "  +7.08%  "
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "'0.0760"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").Value = "'4.43"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").Value = "'128.14"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("D40").Value = "'2.88"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E42").Value = "  -4.59%  "
$ws.Range("D43").Value = "'20.52"
$ws.Range("E43").Value = "  -3.79%  "
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").Value = "'1.938.91"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("D47").Value = "'2.94"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'9.13"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.79"
$ws.Range("E49").Value = "  +8.90%  "
$ws.Range("D50").Value = "'76.77"
$ws.Range("E50").Value = "  +4.01%  "
$ws.Range("D51").Value = "'4.81"
$ws.Range("E51").Value = "  +5.39%  "
